$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "L6"

# Header row (new shared strings must land first, in B,C,D,E order)
$ws.Range("B1").Value = 'Form'
$ws.Range("C1").Value = 'Goals scored'
$ws.Range("D1").Value = 'Goals conceded'
$ws.Range("E1").Value = 'Total Goals'

# Column B: Form (one new shared string per team, in team order)
$colB = @('Ajaccio,W D W L D D', 'Amiens,L L W D D L', 'Auxerre,D W W L W D', 'Caen,D L L D L D', 'Chambly,D L W L D W', 'Chateauroux,D L L D L L', 'Clermont,L W D W W W', 'Dunkerque,W D D L D W', 'Grenoble,D W D W L L', 'Guingamp,L D W W D W', 'Le Havre,L D D D D L', 'Nancy,L D L W D L', 'Niort,D L L L D D', 'Paris FC,W D D L D W', 'Pau FC,W D W D W D', 'Rodez,D L D D D D', 'Sochaux,L D D L W L', 'Toulouse,W L W D W L', 'Troyes,W D W W W W', 'Valenciennes,L W L L D W')
for ($i = 0; $i -lt 20; $i++) { $ws.Cells.Item($i + 2, 2).Value = $colB[$i] }

# Column C: Goals scored
$colC = @('Ajaccio,3 1 3 0 2 0', 'Amiens,0 0 3 2 1 0', 'Auxerre,1 4 3 0 2 0', 'Caen,1 1 0 1 0 0', 'Chambly,2 1 2 0 1 4', 'Chateauroux,2 1 0 1 1 1', 'Clermont,1 3 1 1 2 2', 'Dunkerque,2 1 1 1 1 1', 'Grenoble,2 3 1 2 1 1', 'Guingamp,0 1 1 1 0 3', 'Le Havre,0 1 1 1 1 2', 'Nancy,1 0 2 2 1 1', 'Niort,0 0 0 0 1 1', 'Paris FC,2 1 2 0 0 3', 'Pau FC,3 1 2 1 3 0', 'Rodez,2 0 2 1 1 1', 'Sochaux,1 0 1 0 1 0', 'Toulouse,4 0 4 1 1 1', 'Troyes,1 1 1 3 3 2', 'Valenciennes,0 2 0 1 1 3')
for ($i = 0; $i -lt 20; $i++) { $ws.Cells.Item($i + 2, 3).Value = $colC[$i] }

# Column D: Goals conceded
$colD = @('Ajaccio,0 1 0 2 2 0', 'Amiens,3 2 1 2 1 3', 'Auxerre,1 0 2 3 1 0', 'Caen,1 3 1 1 1 0', 'Chambly,2 2 0 1 1 2', 'Chateauroux,2 2 1 1 2 2', 'Clermont,2 0 1 0 1 1', 'Dunkerque,1 1 1 2 1 0', 'Grenoble,2 1 1 0 3 2', 'Guingamp,1 1 0 0 0 0', 'Le Havre,2 1 1 1 1 4', 'Nancy,4 0 3 0 1 3', 'Niort,0 4 3 3 1 1', 'Paris FC,1 1 2 1 0 1', 'Pau FC,1 1 1 1 0 0', 'Rodez,2 1 2 1 1 1', 'Sochaux,2 0 1 2 0 1', 'Toulouse,0 1 1 1 0 3', 'Troyes,0 1 0 0 1 1', 'Valenciennes,3 1 1 3 1 1')
for ($i = 0; $i -lt 20; $i++) { $ws.Cells.Item($i + 2, 4).Value = $colD[$i] }

# Column E: Total Goals
$colE = @('Ajaccio,3 2 3 2 4 0', 'Amiens,3 2 4 4 2 3', 'Auxerre,2 4 5 3 3 0', 'Caen,2 4 1 2 1 0', 'Chambly,4 3 2 1 2 6', 'Chateauroux,4 3 1 2 3 3', 'Clermont,3 3 2 1 3 3', 'Dunkerque,3 2 2 3 2 1', 'Grenoble,4 4 2 2 4 3', 'Guingamp,1 2 1 1 0 3', 'Le Havre,2 2 2 2 2 6', 'Nancy,5 0 5 2 2 4', 'Niort,0 4 3 3 2 2', 'Paris FC,3 2 4 1 0 4', 'Pau FC,4 2 3 2 3 0', 'Rodez,4 1 4 2 2 2', 'Sochaux,3 0 2 2 1 1', 'Toulouse,4 1 5 2 1 4', 'Troyes,1 2 1 3 4 3', 'Valenciennes,3 3 1 4 2 4')
for ($i = 0; $i -lt 20; $i++) { $ws.Cells.Item($i + 2, 5).Value = $colE[$i] }

# Column A: rank numbers 1-20, stored as text (re-uses existing shared strings "1".."20")
for ($i = 0; $i -lt 20; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).NumberFormat = "@"
  $ws.Cells.Item($r, 1).Value = [string]($i + 1)
}

# Restore the originally active sheet/tab (unchanged by this edit)
$wb.Worksheets.Item(1).Activate()

Write-Output "L6 sheet added"
